# Add the new "2022-Q3" sheet with its fund-holding data, and update the
# "总计" (total) summary sheet with a new leading row for 2022-Q3.

$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
    $range.Borders.Weight = 2            # xlThin
}

# Force text storage (preserves leading zeros / avoids numeric coercion) for
# columns that the source data stores as plain text, e.g. "270041", "4.37".
function Set-TextFormat($range) {
    $range.NumberFormat = "@"
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before the existing "2022-Q2"
#    sheet (tab order becomes: 总计, 2022-Q3, 2022-Q2, 2021-Q4, 2021-Q3,
#    2021-Q2).
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (same layout/style as the other fund-holding sheets).
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"
Set-HeaderStyle $q3Sheet.Range("B1:H1")

# Columns B:G hold plain text in the source workbook (fund code/name keep
# leading zeros, numeric-looking figures are text too) - force text format
# before writing so Excel doesn't coerce them to numbers.
Set-TextFormat $q3Sheet.Range("B2:G5")

# Data rows (values carried over from 2022-Q2 where unchanged, updated
# figures per the 2022-Q3 filing).
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "270041"
$q3Sheet.Range("C2").Value = "广发消费品精选混合A"
$q3Sheet.Range("D2").Value = "4.37"
$q3Sheet.Range("E2").Value = "72.88"
$q3Sheet.Range("F2").Value = "3.27"
$q3Sheet.Range("G2").Value = "0.1429"
$q3Sheet.Range("H2").Value = 9

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "004320"
$q3Sheet.Range("C3").Value = "前海开源沪港深乐享生活灵活配置混合"
$q3Sheet.Range("D3").Value = "0.27"
$q3Sheet.Range("E3").Value = "71.16"
$q3Sheet.Range("F3").Value = "3.23"
$q3Sheet.Range("G3").Value = "0.0087"
$q3Sheet.Range("H3").Value = 8

$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "001613"
$q3Sheet.Range("C4").Value = "长城久祥灵活配置混合"
$q3Sheet.Range("D4").Value = "0.24"
$q3Sheet.Range("E4").Value = "84.63"
$q3Sheet.Range("F4").Value = "3.23"
$q3Sheet.Range("G4").Value = "0.0078"
$q3Sheet.Range("H4").Value = 8

$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "010022"
$q3Sheet.Range("C5").Value = "广发消费品精选混合C"
$q3Sheet.Range("D5").Value = "0.14"
$q3Sheet.Range("E5").Value = "72.88"
$q3Sheet.Range("F5").Value = "3.27"
$q3Sheet.Range("G5").Value = "0.0046"
$q3Sheet.Range("H5").Value = 9

Set-HeaderStyle $q3Sheet.Range("A2:A5")

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for 2022-Q3 above the old
#    2022-Q2 row, pushing the existing rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2:2").Insert()

# Excel's row-insert copies the format of the row above into the new row;
# clear that back out for the non-index columns (only column A carries the
# bold/centered "index" style in this sheet).
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.16
Set-HeaderStyle $totalSheet.Range("A2")

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.24

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 9
$totalSheet.Range("D4").Value = 0.52

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q3"
$totalSheet.Range("C5").Value = 7
$totalSheet.Range("D5").Value = 0.5600000000000001

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q2"
$totalSheet.Range("C6").Value = 20
$totalSheet.Range("D6").Value = 2.65
